# CryCompanywiseStockReport_1 — stock-quantity correction pass.
#
# A batch of line items had their issued/closing Quantity (col F) reduced
# (stock recount), which changes each line's Value (col G = Rate(D) * Qty(F),
# rounded to 2dp). Every "Sub Total:" row (col B) is the sum of the G values
# for its company's block, and the sheet's two grand-total rows (B741/B742)
# are the sum of all Sub Totals — both were re-derived to stay consistent
# with the corrected line values. Two pairs of rows (152/153 "COL-Colgate
# Zigzag..." and 258/259 "HUL-Bru Inst Poly 50g") had their entire data rows
# swapped between the two sibling lines.
#
# Values below are the exact post-edit numbers, written directly to avoid
# any floating point drift versus the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- ASQUARE FOOD BEVERAGES (ZOFF spices) ---
$ws.Range("F28").Value  = 182
$ws.Range("G28").Value  = 4671.94
$ws.Range("F34").Value  = 11
$ws.Range("G34").Value  = 408.98
$ws.Range("F47").Value  = 14
$ws.Range("G47").Value  = 494.48
$ws.Range("B54").Value  = 115644.07   # Sub Total: ASQUARE FOOD BEVERAGES

# --- BHAWAR SALES CORPORATION / BLUE STAR block ---
$ws.Range("F97").Value  = 197
$ws.Range("G97").Value  = 12548.9
$ws.Range("F115").Value = 64
$ws.Range("G115").Value = 3027.2
$ws.Range("B116").Value = 167482.25   # Sub Total: BHAWAR SALES CORPORATION

# --- Bright Lifecare block: row swap (product codes only) ---
$ws.Range("B152").Value = 53925       # COL-Colgate Zigzag ... (item 120)
$ws.Range("B153").Value = 57756       # COL-Colgate Zigzag ... (item 121)

# --- GHP household products ---
$ws.Range("F195").Value = 47
$ws.Range("G195").Value = 1477.21
$ws.Range("F198").Value = 54
$ws.Range("G198").Value = 2176.2
$ws.Range("F203").Value = 206
$ws.Range("G203").Value = 9576.940000000001
$ws.Range("F204").Value = 43
$ws.Range("G204").Value = 3231.02
$ws.Range("F207").Value = 117
$ws.Range("G207").Value = 2811.51
$ws.Range("B209").Value = 55932.99    # Sub Total

# --- GODREJ ---
$ws.Range("F211").Value = 234
$ws.Range("G211").Value = 26746.2
$ws.Range("F213").Value = 44
$ws.Range("G213").Value = 3032.04
$ws.Range("B217").Value = 80747.34    # Sub Total

# --- HINDUSTAN UNILEVER: full row swap between the two HUL-Bru lines ---
$ws.Range("B258").Value = 61610
$ws.Range("D258").Value = 102.71
$ws.Range("E258").Value = 122.71
$ws.Range("F258").Value = 62
$ws.Range("G258").Value = 6368.02
$ws.Range("B259").Value = 57077
$ws.Range("D259").Value = 93.08
$ws.Range("E259").Value = 111.2
$ws.Range("F259").Value = 1
$ws.Range("G259").Value = 93.08

# --- INFLUENCE (LG appliances) ---
$ws.Range("F324").Value = 2
$ws.Range("G324").Value = 50642.62
$ws.Range("B327").Value = 117006.17   # Sub Total

# --- Nutrigreen / CHUKDE spices block ---
$ws.Range("F351").Value = 461
$ws.Range("G351").Value = 64811.99
$ws.Range("B353").Value = 81866.73    # Sub Total
$ws.Range("F398").Value = 130
$ws.Range("G398").Value = 4864.6
$ws.Range("F399").Value = 182
$ws.Range("G399").Value = 9272.9
$ws.Range("B409").Value = 114371.25   # Sub Total

# --- CREMICA biscuits/cookies ---
$ws.Range("F418").Value = 404
$ws.Range("G418").Value = 5312.6
$ws.Range("F420").Value = 219
$ws.Range("G420").Value = 5759.7
$ws.Range("F422").Value = 262
$ws.Range("G422").Value = 4304.66
$ws.Range("F432").Value = 291
$ws.Range("G432").Value = 7653.3
$ws.Range("F433").Value = 220
$ws.Range("G433").Value = 3614.6
$ws.Range("F434").Value = 468
$ws.Range("G434").Value = 6893.64
$ws.Range("B435").Value = 80387.53    # Sub Total

# --- Octavius rice ---
$ws.Range("F475").Value = 62
$ws.Range("G475").Value = 9126.4
$ws.Range("B479").Value = 43728.64    # Sub Total

# --- UNB McVities ---
$ws.Range("F487").Value = 108
$ws.Range("G487").Value = 2849.04
$ws.Range("B489").Value = 29651.8     # Sub Total

# --- RAJ cleaning supplies ---
$ws.Range("F510").Value = 120
$ws.Range("G510").Value = 2510.4
$ws.Range("F516").Value = 31
$ws.Range("G516").Value = 5942.08
$ws.Range("F517").Value = 165
$ws.Range("G517").Value = 4374.15
$ws.Range("F524").Value = 80
$ws.Range("G524").Value = 2151.2
$ws.Range("B526").Value = 67734.8     # Sub Total

# --- RECKITT (Moov, Dettol, Harpic, Lizol) ---
$ws.Range("F555").Value = 170
$ws.Range("G555").Value = 20158.6
$ws.Range("F559").Value = 111
$ws.Range("G559").Value = 6578.97
$ws.Range("F563").Value = 207
$ws.Range("G563").Value = 11360.16
$ws.Range("F564").Value = 116
$ws.Range("G564").Value = 3178.4
$ws.Range("F567").Value = 47
$ws.Range("G567").Value = 6266.51
$ws.Range("F571").Value = 36
$ws.Range("G571").Value = 5103
$ws.Range("B576").Value = 153116.94   # Sub Total

# --- TATA CONSUMER PRODUCT (dals, salt, chilli powder) ---
$ws.Range("F646").Value = 56
$ws.Range("G646").Value = 4567.36
$ws.Range("F649").Value = 66
$ws.Range("G649").Value = 8613
$ws.Range("F653").Value = 33
$ws.Range("G653").Value = 2992.44
$ws.Range("F654").Value = 173
$ws.Range("G654").Value = 12033.88
$ws.Range("F657").Value = 119
$ws.Range("G657").Value = 16066.19
$ws.Range("F659").Value = 15
$ws.Range("G659").Value = 1810.65
$ws.Range("B660").Value = 98059.42999999999  # Sub Total

# --- Tip Top Food Tech (Orgfeed / Shankys) ---
$ws.Range("F665").Value = 64
$ws.Range("G665").Value = 6963.84
$ws.Range("F681").Value = 77
$ws.Range("G681").Value = 3146.99
$ws.Range("F683").Value = 39
$ws.Range("G683").Value = 684.0599999999999
$ws.Range("B690").Value = 80048.71000000001  # Sub Total

# --- VVD AND SONS (edible oils) ---
$ws.Range("F728").Value = 2310
$ws.Range("G728").Value = 376784.1
$ws.Range("F729").Value = 268
$ws.Range("G729").Value = 75809.16
$ws.Range("F730").Value = 366
$ws.Range("G730").Value = 52941.9
$ws.Range("F732").Value = 136
$ws.Range("G732").Value = 10491.04
$ws.Range("B736").Value = 545693.26   # Sub Total

# --- Sheet-wide totals ---
$ws.Range("B741").Value = 3170066.91  # Sub Total (grand, all companies)
$ws.Range("B742").Value = 3170066.91  # Grand Total:
